$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows above row 12 ("Programa resumido:" block) to make room for
# a new "Docentes responsáveis:" section listing the responsible teachers.
$ws.Range("A12:A16").EntireRow.Insert()

# Row 12: section label in column A (inherits the bold style used by the
# other section labels from the row above).
$ws.Range("A12").Value = "Docentes responsáveis:"

# Column A is unused on rows 13-16; make sure no stray cell is left behind.
$ws.Range("A13:A16").Clear()

# Rows 13-16: one teacher per row, duplicated into column B and C (current /
# modified ementa columns).
$ws.Range("B13").Value = "471420 - Carlos Antonio Reis Pereira Baptista"
$ws.Range("C13").Value = "471420 - Carlos Antonio Reis Pereira Baptista"

$ws.Range("B14").Value = "3480026 - João Paulo Pascon"
$ws.Range("C14").Value = "3480026 - João Paulo Pascon"

$ws.Range("B15").Value = "5840793 - Sérgio Schneider"
$ws.Range("C15").Value = "5840793 - Sérgio Schneider"

$ws.Range("B16").Value = "7797767 - Viktor Pastoukhov"
$ws.Range("C16").Value = "7797767 - Viktor Pastoukhov"

# Apply the same cell formatting used throughout the sheet: column B
# wrap-text (style index 2) and column C wrap-text red font (style index 3).
$ws.Range("B3").Copy()
$ws.Range("B13:B16").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C13:C16").PasteSpecial(-4122)
$excel.CutCopyMode = $false
